$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "So tien ghi co: {Amount_}    {Currency_}"  (4 spaces -> 1 space)
#    Only the whitespace-only run's text is touched; {Amount_} / {Currency_}
#    stay in their own runs.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$amountIdx1 = $full.IndexOf("{Amount_}")
$gapStart1 = $amountIdx1 + 9
$gapRange1 = $d.Range($gapStart1, $gapStart1 + 4)
$gapRange1.Text = " "

# ---------------------------------------------------------------------------
# 2) "So tien thuc nop: {Amount_}   {Currency_}" (3 spaces -> 2 spaces)
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$amountIdx2 = $full.IndexOf("Amount_}   {Currency_}")
$gapStart2 = $amountIdx2 + 8
$gapRange2 = $d.Range($gapStart2, $gapStart2 + 3)
$gapRange2.Text = "  "

# ---------------------------------------------------------------------------
# 3) Address line: "{Street _}, {Towndist_}, {City_}, {Country_}"
#       -> "{Street}, {Towndist_}, {City_}, {Country_}"
#    The placeholder "{Street _}" (with the stray trailing " _" before the
#    closing brace) becomes "{Street}, " while the "_GoBack" bookmark that
#    used to wrap "Street " collapses to an empty bookmark right after the
#    new comma, i.e. immediately before "{Towndist_}". Done last since
#    adding the bookmark nudges later same-format run merging.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$streetIdx = $full.IndexOf("{Street")
$streetRange = $d.Range($streetIdx, $streetIdx + 12)   # "{Street _}, "
$streetRange.Text = "{Street}, "

$full = $d.Content.Text
$townIdx = $full.IndexOf("{Towndist_}")
$goBackPoint = $d.Range($townIdx, $townIdx)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
